# Contoso Chai Tea market trends 2023 — header row text refresh
#
# The "Data" header (A1) is left untouched; the other five table headers
# (B1:F1) are reworded / recapitalised per the source commit:
#   - "Total de vendas de chai (unidades)"            -> "Total de vendas de Chai (unidades)"
#   - "Vendas de chai artesanal (unidades)"            -> "Vendas de Chai Artesanal (unidades)"
#   - "Vendas de chai pronto (unidades)"                -> "Vendas de Chai pré-fabricado (unidades)"
#   - "Participação nas redes sociais (exibições)"      -> "Engajamento em redes sociais (visualizações)"
#   - "Pesquisas online por chai"                       -> "Pesquisas online para Chai"
#
# These header cells are the header row of Table1, so updating the cell
# values also renames the corresponding ListColumns / table XML in one go.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Total de vendas de Chai (unidades)"
$ws.Range("C1").Value = "Vendas de Chai Artesanal (unidades)"
$ws.Range("D1").Value = "Vendas de Chai pré-fabricado (unidades)"
$ws.Range("E1").Value = "Engajamento em redes sociais (visualizações)"
$ws.Range("F1").Value = "Pesquisas online para Chai"
